$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new "Exceptions handling" row to the Tehnical functionality table
$ws.Range("A40").Value = 10
$ws.Range("B40").Value = "Exceptions handling for better FE control of errors"
$ws.Range("B40").WrapText = $true

# Match the formatting (fill + medium border) used by the other rows in the table,
# copying it from the row directly above (A39:C39)
$ws.Range("C39").Copy()
$ws.Range("C40").PasteSpecial(-4122)

# Update the active selection to reflect where the user ended up editing
$ws.Range("D21").Select() | Out-Null
